$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15; this shifts rows 15..60 down to 16..61
$ws.Rows.Item(15).Insert()

# Fill the new row 15 with data
$ws.Cells.Item(15, 1).Value = 4
$ws.Cells.Item(15, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(15, 3).Value = "Los Lagos"
$ws.Cells.Item(15, 4).Value = 44707
$ws.Cells.Item(15, 5).Value = 10
$ws.Cells.Item(15, 6).Value = 100112043
$ws.Cells.Item(15, 7).Value = "Pepino dulce"
$ws.Cells.Item(15, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(15, 9).Value = "Especial"
$ws.Cells.Item(15, 10).Value = 50
$ws.Cells.Item(15, 11).Value = 21000
$ws.Cells.Item(15, 12).Value = 21000
$ws.Cells.Item(15, 13).Value = 21000
$ws.Cells.Item(15, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 16).Value = 1167
$ws.Cells.Item(15, 17).Value = 18
$ws.Cells.Item(15, 18).Value = "Hortaliza"
